$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like numbers need to be forced to Text
# format before assignment, otherwise Excel auto-converts them to numeric
# values (losing the original formatted text, e.g. trailing zeros).

$ws.Range("D2").Value = '30.700.08'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '2.123.36'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.015'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.012'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5265'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4563'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '55.09'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09134'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.177'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.53'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").Value = '2.128.19'
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.873'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.165'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001179'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '97.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.014'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06711'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.011'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.347'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("D23").Value = '30.774.99'
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.368'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.40%  '
$ws.Range("D26").Value = '2.362.26'
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '165.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.568'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.211'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1077'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.658'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.389'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.950'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.912'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02666'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06891'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2333'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6935'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.261'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.65%  '
$ws.Range("E45").Value = '  +2.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.320'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000378'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +22.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.703'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.258'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07320'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.91%  '
